$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (llama-3.1-8b-instant)
$ws.Range("I2").Value = 15
$ws.Range("J2").Value = 0.001041666666666667
$ws.Range("K2").Value = 4400
$ws.Range("L2").Value = 0.008800000000000001

# Row 4 (meta-llama/llama-4-maverick-17b-128e-instruct)
$ws.Range("I4").Value = 2
$ws.Range("J4").Value = 0.002
$ws.Range("K4").Value = 540
$ws.Range("L4").Value = 0.00108

# Row 14 (qwen/qwen3-32b)
$ws.Range("I14").Value = 2
$ws.Range("J14").Value = 0.002
$ws.Range("K14").Value = 290
$ws.Range("L14").Value = 0.00058

$wb.Save()
